# Generate Report for Handback
# Marks the d9731ca5 file as handed back (in sync with en-US) for both the
# zh-cn and de-de localization sheets: updates Status, fills in the
# "Latest Target File" / "Latest Handback File" columns, stamps a real
# "Latest Handback DateTime", and adds the corresponding hyperlink on the
# newly-populated Latest Target File cell.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"
$targetFileName   = "d9731ca5-0798-4253-8d84-6df230963dc5.md"
$targetUrl        = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e0f417447bf629aa2238499d7ce7c435007a73ff/e2e/d9731ca5-0798-4253-8d84-6df230963dc5.md"

function Update-HandbackSheet($SheetName, $XlfFile, $HandbackDateTime) {
    $ws = $wb.Worksheets.Item($SheetName)

    # --- Row 4 (d9731ca5-...-dc5.md) ---
    $ws.Range("C4").Value = $statusHandedBack
    $ws.Range("J4").Value = $XlfFile
    $ws.Range("K4").Value = $HandbackDateTime

    # --- Row 5 (f416fa0e-...-432.md, depends on the same d9731ca5 xlf) ---
    $ws.Range("C5").Value = $statusHandedBack
    $ws.Range("J5").Value = $XlfFile
    $ws.Range("K5").Value = $HandbackDateTime

    # The "Latest Target File" column (I) now references the markdown file,
    # and needs the hyperlink that goes with it (same as column A's link).
    $ws.Hyperlinks.Add($ws.Range("I4"), $targetUrl, [Type]::Missing, [Type]::Missing, $targetFileName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I5"), $targetUrl, [Type]::Missing, [Type]::Missing, $targetFileName) | Out-Null
}

Update-HandbackSheet "zh-cn" "d9731ca5-0798-4253-8d84-6df230963dc5.ce079edaad233f35b74cf10d4be1d31439fe14e5.zh-cn.xlf" "2016-11-09 06:08:13"
Update-HandbackSheet "de-de" "d9731ca5-0798-4253-8d84-6df230963dc5.ce079edaad233f35b74cf10d4be1d31439fe14e5.de-de.xlf" "2016-11-09 06:08:32"
